{"js": "const replacements = [\n  [\"78\u00f79=8, 6\", \"15\u00f79=1, 6\"],\n  [\"97\u00f72=48, 1\", \"72\u00f74=18, 0\"],\n  [\"26\u00f73=8, 2\", \"74\u00f79=8, 2\"],\n  [\"85\u00f78=10, 5\", \"83\u00f78=10, 3\"],\n  [\"11\u00f77=1, 4\", \"77\u00f75=15, 2\"],\n  [\"76\u00f75=15, 1\", \"40\u00f76=6, 4\"],\n  [\"10\u00f78=1, 2\", \"21\u00f78=2, 5\"],\n  [\"89\u00f76=14, 5\", \"68\u00f76=11, 2\"],\n  [\"26\u00f76=4, 2\", \"46\u00f79=5, 1\"],\n  [\"75\u00f75=15, 0\", \"91\u00f77=13, 0\"],\n  [\"15\u00f76=2, 3\", \"13\u00f79=1, 4\"],\n  [\"27\u00f72=13, 1\", \"99\u00f72=49, 1\"],\n  [\"43\u00f73=14, 1\", \"67\u00f79=7, 4\"],\n  [\"25\u00f73=8, 1\", \"87\u00f72=43, 1\"],\n  [\"55\u00f74=13, 3\", \"83\u00f76=13, 5\"],\n  [\"82\u00f79=9, 1\", \"21\u00f76=3, 3\"],\n  [\"99\u00f79=11, 0\", \"82\u00f75=16, 2\"],\n  [\"14\u00f74=3, 2\", \"64\u00f77=9, 1\"],\n  [\"82\u00f74=20, 2\", \"25\u00f77=3, 4\"],\n  [\"74\u00f77=10, 4\", \"28\u00f72=14, 0\"],\n  [\"45\u00f79=5, 0\", \"74\u00f74=18, 2\"],\n  [\"93\u00f74=23, 1\", \"18\u00f72=9, 0\"],\n  [\"73\u00f78=9, 1\", \"87\u00f77=12, 3\"],\n  [\"51\u00f72=25, 1\", \"36\u00f73=12, 0\"],\n  [\"80\u00f77=11, 3\", \"52\u00f78=6, 4\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    ,@(\"78\u00f79=8, 6\", \"15\u00f79=1, 6\")\n    ,@(\"97\u00f72=48, 1\", \"72\u00f74=18, 0\")\n    ,@(\"26\u00f73=8, 2\", \"74\u00f79=8, 2\")\n    ,@(\"85\u00f78=10, 5\", \"83\u00f78=10, 3\")\n    ,@(\"11\u00f77=1, 4\", \"77\u00f75=15, 2\")\n    ,@(\"76\u00f75=15, 1\", \"40\u00f76=6, 4\")\n    ,@(\"10\u00f78=1, 2\", \"21\u00f78=2, 5\")\n    ,@(\"89\u00f76=14, 5\", \"68\u00f76=11, 2\")\n    ,@(\"26\u00f76=4, 2\", \"46\u00f79=5, 1\")\n    ,@(\"75\u00f75=15, 0\", \"91\u00f77=13, 0\")\n    ,@(\"15\u00f76=2, 3\", \"13\u00f79=1, 4\")\n    ,@(\"27\u00f72=13, 1\", \"99\u00f72=49, 1\")\n    ,@(\"43\u00f73=14, 1\", \"67\u00f79=7, 4\")\n    ,@(\"25\u00f73=8, 1\", \"87\u00f72=43, 1\")\n    ,@(\"55\u00f74=13, 3\", \"83\u00f76=13, 5\")\n    ,@(\"82\u00f79=9, 1\", \"21\u00f76=3, 3\")\n    ,@(\"99\u00f79=11, 0\", \"82\u00f75=16, 2\")\n    ,@(\"14\u00f74=3, 2\", \"64\u00f77=9, 1\")\n    ,@(\"82\u00f74=20, 2\", \"25\u00f77=3, 4\")\n    ,@(\"74\u00f77=10, 4\", \"28\u00f72=14, 0\")\n    ,@(\"45\u00f79=5, 0\", \"74\u00f74=18, 2\")\n    ,@(\"93\u00f74=23, 1\", \"18\u00f72=9, 0\")\n    ,@(\"73\u00f78=9, 1\", \"87\u00f77=12, 3\")\n    ,@(\"51\u00f72=25, 1\", \"36\u00f73=12, 0\")\n    ,@(\"80\u00f77=11, 3\", \"52\u00f78=6, 4\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $found = $find.Execute([ref]$oldText, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$newText, [ref]2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
